# Update column F (dSF) values for several rows as part of a data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 4
    7  = -2
    14 = -1
    20 = 4
    44 = -7
    48 = 0
    49 = -1
    50 = -3
    51 = -1
    54 = 1
    62 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
